$wb = $excel.ActiveWorkbook

# --- Sheet "Issue Tracking": update comment cells per call with Thomas ---
$ws = $wb.Worksheets.Item("Issue Tracking")

# Row 11: add Brandon follow-up comment in the "Resolution / Plan" column (E)
$ws.Range("E11").Value = "•12/22 Brandon to Verify if this is the case."

# Row 13: append "Can we close this item?" to the firewall-rules resolution note (E13)
$ws.Range("E13").Value = "•The comment as before.  The firewall rules will certainly not be correct.`n•12/14/ Not able to reproduce with my device and the configuratio file from the customer.`n•12/22 Can we close this item?"

# Row 12: append "Can we close this item?" to the existing resolution note (E12)
$ws.Range("E12").Value = "•A configuration file made on a different CloudGate type will not work on any CloudGate. This was already the case with configuration files from a CG 3G placed on a CG LTE.`n•12/14/ Not able to reproduce with my device and the configuratio file from the customer.`n•12/15 Sent customer logs for verification with the App team.`n•12/19 Customer to try the latest LuvitRED 2.10.0, but more important to try a configuration file created on a CloudGate Gemalto. `n•12/22 Can we close this item?"

# Row 12 is now taller since the note grew by one line
$ws.Rows.Item(12).RowHeight = 210

# Scroll position / selection moved back up after the call
$ws.Range("F12").Select()
